# [ADD] New normalize way
# Update the normalized values in column B (rows 3-10) to reflect the new
# normalization approach, while leaving all other data/styles untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 0.35902840511623263
$ws.Range("B4").Value = 0.9152830807763157
$ws.Range("B5").Value = 1.8305661615526314
$ws.Range("B6").Value = 2.1541704306973957
$ws.Range("B7").Value = 2.7458492423289473
$ws.Range("B8").Value = 2.9919033759686053
$ws.Range("B9").Value = 3.661132323105263
$ws.Range("B10").Value = 4.576415403881579
